$wb = $excel.ActiveWorkbook

# New row data per sheet (1-based sheet index matches workbook order:
# 1=ROW35-FE-LIFTER, 2=ROW35-MID-LIFTER, 3=ROW02-FE-LIFTER, 4=ROW02-MID-LIFTER)

$rowNum = 56

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item($rowNum, 1).Value = 45752.36009445602
$ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($rowNum, 2).Value = "0x01,0x90"
$ws.Cells.Item($rowNum, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item($rowNum, 4).Value = "0x01,0x6e"
$ws.Cells.Item($rowNum, 5).Value = "0xd"
$ws.Cells.Item($rowNum, 6).Value = 400
$ws.Cells.Item($rowNum, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($rowNum, 8).Value = 366
$ws.Cells.Item($rowNum, 9).Value = 13

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item($rowNum, 1).Value = 45752.21420876157
$ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($rowNum, 2).Value = "0x01,0x90"
$ws.Cells.Item($rowNum, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($rowNum, 4).Value = "0x01,0x6e"
$ws.Cells.Item($rowNum, 5).Value = "0xe"
$ws.Cells.Item($rowNum, 6).Value = 400
$ws.Cells.Item($rowNum, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($rowNum, 8).Value = 366
$ws.Cells.Item($rowNum, 9).Value = 14

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item($rowNum, 1).Value = 45752.35381758102
$ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($rowNum, 2).Value = "0x01,0x90"
$ws.Cells.Item($rowNum, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item($rowNum, 4).Value = "0x01,0x6e"
$ws.Cells.Item($rowNum, 5).Value = "0x3"
$ws.Cells.Item($rowNum, 6).Value = 400
$ws.Cells.Item($rowNum, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($rowNum, 8).Value = 366
$ws.Cells.Item($rowNum, 9).Value = 3

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item($rowNum, 1).Value = 45752.41215091435
$ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($rowNum, 2).Value = "0x01,0x90"
$ws.Cells.Item($rowNum, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item($rowNum, 4).Value = "0x01,0x6e"
$ws.Cells.Item($rowNum, 5).Value = "0x3"
$ws.Cells.Item($rowNum, 6).Value = 400
$ws.Cells.Item($rowNum, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item($rowNum, 8).Value = 366
$ws.Cells.Item($rowNum, 9).Value = 3
